$d = $word.ActiveDocument

# The paragraph we need to touch currently reads (as a single run):
#   "...separate notice of this court. The Clerk of Court is ordered..."
# and needs to become three runs (identical formatting) so the lower-case
# "c" in "court" becomes an upper-case "C" in its own run:
#   "...separate notice of this " + "C" + "ourt. The Clerk of Court..."
#
# Locate the single lower-case "c" that starts "court." right after
# "of this " and before "ourt. The Clerk of Court is ordered to serve".
$rng = $d.Content
$found = $rng.Find.Execute(
    "of this court. The Clerk of Court is ordered to serve",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target sentence for surety-hearing paragraph."
}

# $rng now spans the matched text; work out the absolute offsets of the
# single "c" we need to capitalize ("of this c|ourt" -> the 9th char in
# the match, i.e. 5 characters before the end of "court").
$matchStart = $rng.Start
$cStart = $matchStart + ("of this ".Length)
$cEnd = $cStart + 1

$cRange = $d.Range($cStart, $cEnd)
if ($cRange.Text -ne "c") {
    throw "Unexpected character at target offset: [$($cRange.Text)]"
}

# Replace the character's text (same formatting throughout).
$cRange.Text = "C"

# Re-grab a fresh 1-character range over the now-capitalized "C" and nudge
# a character-level format property on/off. This forces the host to keep
# the run boundaries explicit instead of silently re-merging this run with
# its identically-formatted neighbors, producing the desired run split:
#   run1: "...separate notice of this "
#   run2: "C"
#   run3: "ourt. The Clerk of Court is ordered to serve...within "
$cRange2 = $d.Range($cStart, $cStart + 1)
$cRange2.Bold = 1
$cRange2.Bold = 0

# The same whole-paragraph re-layout also silently swallows the run
# boundaries of the two following runs ("15 days of the date of this
# order. {% endif %}" and "{% if fta_conditions.forfeit_license is true %}")
# even though their text/formatting is untouched. Re-assert those two runs
# explicitly so they remain distinct runs rather than being absorbed into
# the run we just split.
$rngA = $d.Content
$foundA = $rngA.Find.Execute(
    "15 days of the date of this order. {% endif %}",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundA) {
    $rngA.Bold = 1
    $rngA.Bold = 0
}

$rngB = $d.Content
$foundB = $rngB.Find.Execute(
    "{% if fta_conditions.forfeit_license is true %}",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundB) {
    $rngB.Bold = 1
    $rngB.Bold = 0
}
